# Updated cryptos list values (Price / Volume(1h)) per the latest scrape.
# D-column values that look numeric are prefixed with a leading apostrophe
# so Excel stores them as text (matching the workbook's original inline
# string / text-formatted cells) instead of auto-converting them to
# floating point numbers and losing formatting (trailing zeros, thousands
# separators written with dots, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.985.65"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "1.823.41"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  -0.64%  "
$ws.Range("D5").Value = "'311.13"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "'0.4251"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("D8").Value = "'0.3658"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "'0.07227"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").Value = "'0.8411"
$ws.Range("D11").Value = "'20.57"
$ws.Range("E11").Value = "  -3.23%  "
$ws.Range("D12").Value = "1.823.07"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "'0.07050"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "'5.276"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "'89.63"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D18").Value = "'0.000008729"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").Value = "27.057.98"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").Value = "'5.119"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "'10.79"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").Value = "2.056.43"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").Value = "'1.979"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").Value = "'150.77"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("D27").Value = "'2.217"
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("D29").Value = "'5.211"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "'116.71"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").Value = "'0.08708"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("D32").Value = "'1.173"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("D33").Value = "'0.7352"
$ws.Range("E33").Value = "  -4.25%  "
$ws.Range("D34").Value = "'2.901"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'4.410"
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("D38").Value = "'0.01936"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").Value = "'0.05213"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").Value = "'7.220"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").Value = "'2.864"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "'0.1685"
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").Value = "'0.5114"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "'8.522"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("D45").Value = "'10.56"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "'1.953"
$ws.Range("E46").Value = "  +6.71%  "
$ws.Range("D47").Value = "'0.4728"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "'105.69"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").Value = "'0.9998"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").Value = "'0.06315"
$ws.Range("D51").Value = "'1.649"
$ws.Range("E51").Value = "  -1.29%  "
